$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Annulla prestiti" -> "Termina prestiti" (keeps the existing bold
#    run formatting).
# ------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Annulla", $true, $true, $false, $false, $false, $true, 1, $false, "Termina", 2)

# ------------------------------------------------------------------
# 2) Explanatory sentence: replace the whole tail (this also drops the
#    old final period) with the complete new sentence in one go, so
#    the run keeps its existing sz/szCs run formatting.
# ------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("annullare tutti i suoi prestiti.", $true, $true, $false, $false, $false, $true, 1, $false, "annullare tutti i suoi prestiti, facendo tornare le risorse associate in archivio.", 2)

# ------------------------------------------------------------------
# 3) Force the explanatory sentence to split into three runs at the
#    right boundaries by dropping a temporary bookmark at each split
#    point (inserting a bookmark splits the underlying run without
#    touching its formatting) and then removing the temporary
#    bookmark again.
# ------------------------------------------------------------------
$split1 = $d.Content
$split1.Find.Execute("prestiti, facendo", $true, $true)
$split1.Collapse(1)
$split1.MoveStart(1, 8)
$d.Bookmarks.Add("ztmpsplit1", $split1)
$d.Bookmarks("ztmpsplit1").Delete()

$split2 = $d.Content
$split2.Find.Execute("in archivio.", $true, $true)
$split2.Collapse(0)
$split2.MoveStart(1, -1)
$d.Bookmarks.Add("ztmpsplit2", $split2)
$d.Bookmarks("ztmpsplit2").Delete()

# ------------------------------------------------------------------
# 4) Place the _GoBack bookmark right after "Termina" (before
#    " prestiti"). Adding it here also removes the old _GoBack
#    bookmark that used to sit near "diversi", since bookmark names
#    must be unique within the document.
# ------------------------------------------------------------------
$goBackRng = $d.Content
$goBackRng.Find.Execute("Termina", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$goBackRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBackRng)
